# Insert a new weekly price-record row right after row 70 (i.e. as the new
# row 71), pushing the existing rows 71-106 down to 72-107. This mirrors a
# new "Región del Maule" / "$/docena de atados (4 kilos)" sample being added
# to the front of the series.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 71:106 down to 72:107 and open up a blank row 71.
$ws.Rows.Item(71).Insert()

# Populate the newly inserted row 71 with the new record.
$ws.Cells.Item(71, 1).Value  = 4
$ws.Cells.Item(71, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(71, 3).Value  = "Los Lagos"
$ws.Cells.Item(71, 4).Value  = 44466
$ws.Cells.Item(71, 5).Value  = 10
$ws.Cells.Item(71, 6).Value  = 100112009
$ws.Cells.Item(71, 7).Value  = "Acelga"
$ws.Cells.Item(71, 8).Value  = "Sin especificar"
$ws.Cells.Item(71, 9).Value  = "Primera"
$ws.Cells.Item(71, 10).Value = 100
$ws.Cells.Item(71, 11).Value = 4000
$ws.Cells.Item(71, 12).Value = 4000
$ws.Cells.Item(71, 13).Value = 4000
$ws.Cells.Item(71, 14).Value = "`$/docena de atados (4 kilos)"
$ws.Cells.Item(71, 15).Value = "Región del Maule"
$ws.Cells.Item(71, 16).Value = 1000
$ws.Cells.Item(71, 17).Value = 4
$ws.Cells.Item(71, 18).Value = "Hortaliza"
